$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "FALSE_count"
$ws.Range("C1").Value = "FALSE_percent"
$ws.Range("D1").Value = "TRUE_count"
$ws.Range("E1").Value = "TRUE_percent"

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 5.555555555555555
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 6.756756756756757

$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 72.22222222222221
$ws.Range("D3").Value = 52
$ws.Range("E3").Value = 70.27027027027027

$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 22.22222222222222
$ws.Range("D4").Value = 13
$ws.Range("E4").Value = 17.56756756756757

$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 5.405405405405405
